# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record": the old
# scraper only pulled team statistics, not the season W-L-T record, so
# three new columns (AD, AE, AF) are appended with the team's record
# repeated on every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the three new header cells so they match the style
# used by the rest of row 1 (style index 1 in the original workbook).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team's 2003 season record (101 wins, 61 losses, 1 tie) repeated down
# every player row.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 101
    $ws.Cells.Item($row, 31).Value = 61
    $ws.Cells.Item($row, 32).Value = 1
}
